$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.202.68"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "3.146.01"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  +0.08%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "592.42"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.22%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "138.28"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.54%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.140.56"
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("E10").Value = "  -1.39%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.28"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("E13").Value = "  -2.34%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "34.22"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").Value = "3.665.45"
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("D17").Value = "3.140.60"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "63.168.90"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("E19").Value = "  -1.90%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "473.39"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.80%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "14.09"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("E22").Value = "  -0.71%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "7.68"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.08%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "84.80"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.97%  "
$ws.Range("E25").Value = "  -1.51%  "
$ws.Range("E26").Value = "  +0.01%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.72"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "8.00"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -2.75%  "
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "6.99"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("E30").Value = "  +3.84%  "
$ws.Range("E31").Value = "  -0.01%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "26.85"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.00%  "
$ws.Range("E33").Value = "  -4.12%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "2.54"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -3.56%  "
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("E36").Value = "  -2.56%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "52.35"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("D38").Value = "0.0₃0697"
$ws.Range("E38").Value = "  -6.46%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0388"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.19%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "420.86"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -3.10%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.76"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -5.98%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "8.22"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.57%  "
$ws.Range("D43").Value = "2.913.67"
$ws.Range("E43").Value = "  +1.73%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.112"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -4.97%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.263"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.46%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.13"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -2.71%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.02%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "25.46"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.84%  "
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("E50").Value = "  -6.70%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "120.42"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.53%  "
